$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -2.611543900785094
$ws.Range("D2").Value = 0.07717835705609512
$ws.Range("E2").Value = 0.9990623130852651

$ws.Range("C3").Value = -0.3272621072280406
$ws.Range("D3").Value = 0.4833682071129712
$ws.Range("E3").Value = -0.5873498981461434

$ws.Range("C4").Value = -112468.1821454628
$ws.Range("D4").Value = 0.9996504337695269
$ws.Range("E4").Value = -0.01202488522617391

$ws.Range("C5").Value = -24.58953277079197
$ws.Range("D5").Value = 0.7758270965853896
$ws.Range("E5").Value = -0.5637216189927576

$ws.Range("C6").Value = -4.042839460422329
$ws.Range("D6").Value = 0.594249163167603
$ws.Range("E6").Value = 0.4776221073808626

$ws.Range("C7").Value = -7010.93806296623
$ws.Range("D7").Value = 0.8004818402675289
$ws.Range("E7").Value = -0.9643514049377471

$ws.Range("C8").Value = -0.002379704966644303
$ws.Range("D8").Value = 0.885039846324471
$ws.Range("E8").Value = -0.3662337739607447

$ws.Range("C9").Value = -4539.519187915086
$ws.Range("D9").Value = 0.5098657370392975
$ws.Range("E9").Value = 0.9605638732232217
